$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set value "done" in C12
$ws.Range("C12").Value = "done"

# Move the visible top-left cell and selection as recorded in the diff
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("C13").Select()
